$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Step 1: update the inputs on row 5 (Contador=1, Tamanho=5); dependent formulas recalc automatically
$ws.Range("E5").Value = 1
$ws.Range("H5").Value = 5

# Step 2: delete the now-duplicate row 6 entirely, shifting rows 7:28 up to 6:27
$ws.Rows.Item(6).Delete()

# Step 3: restore the literal text in the "Primeira parte" / "Arquivo-fonte" columns
# for the rows that shifted up (the runtime does not re-resolve the relative
# "value of the row above" formulas correctly after the row delete)
for ($r = 6; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "cut -c1-"
    $ws.Cells.Item($r, 4).Value = " < arquivo_referencia/pi-1M.txt > "
}

# Force a full recalculation so dependent formulas (J column) drop any
# stale #REF! cached from the row delete above.
$excel.CalculateFull()

# Step 4: restore the selection
$ws.Range("B15").Select()

$wb.Save()
